# Update computed loading_percent results for the 380 kV case (Case_1_230)
# Columns updated per row: B, D, E, F, G, J, L, N, O (C, H, I, K, M stay 0; A is the index column)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B", "D", "E", "F", "G", "J", "L", "N", "O")

# Hashtable keyed by sheet row number -> new values for $columns, in order
$newValues = @{
    2 = @(18.69052045650762, 9.78185833145368, 15.37272156341437, 33.08862196407312, 3.643637610975123, 11.10974437661918, 10.91022497406259, 18.82105546858823, 24.9151233548535)
    3 = @(18.43430649547129, 9.794295355105081, 15.37402584538306, 33.10013831192594, 3.64623275787446, 11.11235173796157, 10.71875143154843, 18.84925070062296, 24.89717123379022)
    4 = @(18.27801891136128, 9.802526863294304, 15.37581688423965, 33.11638535905057, 3.647912149342831, 11.11455741830461, 10.6004293971159, 18.86879421271721, 24.89288586044524)
    5 = @(18.21466068635279, 9.806031285499316, 15.37679709219353, 33.12531033098262, 3.648618201994913, 11.11560904730148, 10.55207803656423, 18.87732049113218, 24.89283464529509)
    6 = @(18.20416213271962, 9.806622263505705, 15.37697501419007, 33.12693136001691, 3.648736753443352, 11.11579291910952, 10.54404288721049, 18.87877025385351, 24.89292852306718)
    7 = @(18.2771630087646, 9.802573517249195, 15.37582908825207, 33.11649640020968, 3.647921583511941, 11.11457098131693, 10.59977778210158, 18.86890692371, 24.89287830609027)
    8 = @(18.60200728900038, 9.786023319167425, 15.37296657840088, 33.09068682916628, 3.64451461953177, 11.11051831672231, 10.84439515470271, 18.83031463922149, 24.90753611011766)
    9 = @(19.2440522863444, 9.75827422859083, 15.37515074120853, 33.11296295737508, 3.63851238060577, 11.10733775656572, 11.31569365671576, 18.77229740638588, 24.98963757735438)
    10 = @(19.71466812311204, 9.740733682214438, 15.38142940715087, 33.17378474925947, 3.634511811515067, 11.10786461600533, 11.65375535510399, 18.74038024417421, 25.08226700299093)
    11 = @(19.92766688003119, 9.733367505365297, 15.38528318924443, 33.21107930885669, 3.632779753137218, 11.10871693483325, 11.80515549335509, 18.72817238924528, 25.13134645251834)
    12 = @(20.00809827319289, 9.730665920756898, 15.38688442491262, 33.22658093434428, 3.632136422564101, 11.10912698153891, 11.86209579518593, 18.72388081860122, 25.15092105378349)
    13 = @(19.99078710805588, 9.731243854578979, 15.38653328463829, 33.22318115268625, 3.632274417648044, 11.10903480126886, 11.8498508466012, 18.72479036929783, 25.14666147208161)
    14 = @(19.93428895567378, 9.733143486120397, 15.38541208918039, 33.21232703606757, 3.632726574552621, 11.10874892481877, 11.80984812437031, 18.72781268624901, 25.13293709635976)
    15 = @(19.89965064508526, 9.734318492942197, 15.38474375631184, 33.20585797612248, 3.633005167745741, 11.10858516111755, 11.78529290102325, 18.72970705169926, 25.12465906373698)
    16 = @(19.70072022186836, 9.74122738759284, 15.38119749478643, 33.17154075607965, 3.634626767112914, 11.10782117438627, 11.64380878084645, 18.74122447464043, 25.0791985002446)
    17 = @(19.57835376990933, 9.745622561334304, 15.37927638390223, 33.15295057171026, 3.635644011535038, 11.10750883745862, 11.55636642199453, 18.74888127919852, 25.05308197035)
    18 = @(19.50787361105377, 9.748208277529802, 15.37826540446192, 33.14316449370542, 3.636237373631775, 11.10738690380158, 11.50585007018162, 18.7535029006307, 25.03871445927373)
    19 = @(19.48399555149814, 9.749093681819474, 15.37793929454065, 33.1400069230507, 3.636439698247918, 11.10735555600121, 11.48870952552934, 18.75510511337101, 25.03396244634661)
    20 = @(19.59139056876307, 9.745148715361989, 15.37947117275112, 33.15483574179868, 3.635534868695304, 11.10753611884896, 11.56569812785065, 18.74804368325102, 25.05579449064688)
    21 = @(19.95089053212854, 9.732583137377905, 15.38573757313715, 33.21547777816394, 3.632593424695771, 11.10883053106303, 11.82160889689477, 18.72691597662455, 25.13694150770486)
    22 = @(20.18449269115672, 9.724882571674714, 15.39065933824172, 33.26314503320802, 3.630744213609582, 11.11018492681403, 11.98655789169028, 18.71503821133511, 25.19573691848374)
    23 = @(20.05996115704881, 9.728945792855846, 15.3879574152683, 33.236971161281, 3.631724496969514, 11.10941580898575, 11.89874764918817, 18.72120133810091, 25.1638328725152)
    24 = @(19.58549703165879, 9.745362757741796, 15.37938281739427, 33.15398064732853, 3.635584185582557, 11.10752360535559, 11.56148002170325, 18.74842167617549, 25.05456614342147)
    25 = @(19.07025732658012, 9.765279591684116, 15.37373319047165, 33.09912510719557, 3.640063945741422, 11.10769219638327, 11.18941972295424, 18.72479036929783, 25.14666147208161)
}

foreach ($row in $newValues.Keys) {
    $values = $newValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}

Write-Output "Updated loading_percent values for rows 2-25 (380 kV case)"
